# Working on state machine: add a new "in/out/factor" block (rows 13-14)
# mirroring the header row, with a new divider computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: repeat header labels (same shared strings as row 1)
$ws.Range("A13").Value = "in "
$ws.Range("B13").Value = "out"
$ws.Range("C13").Value = "factor"

# Row 14: new data row with a divider formula like the rows above
$ws.Range("A14").Value = 50000000
$ws.Range("B14").Value = 22000
$ws.Range("C14").Formula = "=A14/(B14*2)"

# Match the number formatting used by the other "factor" column cells (style index 1 -> numFmtId 1)
$ws.Range("C14").NumberFormat = "0"

# Update the selection to match the newly edited range
$ws.Range("A13:C14").Select()

$wb.Save()
